$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Row 12 (col J): average of the k-column (J2:J11), bold
# ------------------------------------------------------------------
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# Build a "bold" font on a scratch cell and paste that *format only*
# onto J12 in a single atomic step - avoids leaving unused/orphaned
# style records behind when several font properties are changed in
# sequence on the same range.
$scratch1 = $ws.Range("ZZ1")
$scratch1.Font.Bold = $true
$scratch1.Copy() | Out-Null
$ws.Range("J12").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$scratch1.Clear() | Out-Null

# ------------------------------------------------------------------
# Summary rows 14-17
# ------------------------------------------------------------------
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Build a bold, size-12, vertically-centered font on a scratch cell
# and paste that format onto B14:B17 in one shot (same reasoning as
# above - keeps the generated cellXfs/font table minimal & clean).
$scratch2 = $ws.Range("ZZ2")
$scratch2.Font.Bold = $true
$scratch2.Font.Size = 12
$scratch2.VerticalAlignment = [Microsoft.Office.Interop.Excel.XlVAlign]::xlVAlignCenter
$scratch2.Copy() | Out-Null
$ws.Range("B14:B17").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$scratch2.Clear() | Out-Null

$ws.Application.CutCopyMode = $false

# Slightly taller rows for the bigger (size-12) summary font
$ws.Rows.Item(14).RowHeight = 15.6
$ws.Rows.Item(15).RowHeight = 15.6
$ws.Rows.Item(16).RowHeight = 15.6
$ws.Rows.Item(17).RowHeight = 15.6

# ------------------------------------------------------------------
# Page setup (paper size / orientation)
# ------------------------------------------------------------------
$ws.PageSetup.PaperSize = [Microsoft.Office.Interop.Excel.XlPaperSize]::xlPaperA4
$ws.PageSetup.Orientation = [Microsoft.Office.Interop.Excel.XlPageOrientation]::xlPortrait

# Leave selection on J12, matching the authored workbook
$ws.Range("J12").Select() | Out-Null
